$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rephrase the six "Comment" category strings used in column E.
# Order matters: replacing in this sequence reproduces the same
# shared-string table ordering as the authored workbook (old,
# now-unreferenced strings are dropped and replacements are appended
# at the end in the order they are introduced).
$ws.Cells.Replace("Network added to Addon Package within Service in Apr 2020", "Network Added to Add-On Package")
$ws.Cells.Replace("Network moved from Base Service to Addon Package in Apr 2020", "Network Moved from Base Service to Add-On Package")
$ws.Cells.Replace("Network added to base Service in Apr 2020", "Network Added to Base Service")
$ws.Cells.Replace("Old Network removed from database in Apr 2020", "Network Removed from Database")
$ws.Cells.Replace("New Alias added for Network in Apr 2020", "Alias Changed for Network")
$ws.Cells.Replace("New Network added to database in Apr 2020", "New Network Added to Database in Apr 2020")

# Update the active selection to E2, matching the saved view state.
$ws.Range("E2").Select()
